# Text updates as supplied by PM&C.
# Update the "Description" sheet's source/reference footer:
#   - A11/B11 used to hold a single "Sourced from ABS Causes of Death, Australia."
#     note. It is split into a labelled "Source" row whose text is refreshed
#     to "ABS (unpublished) Causes of Death, Australia.".
#   - A new "References" row (row 12) is added, citing the COAG 2011
#     National Healthcare Agreement.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

$ws.Range("A11").Value = "Source"
$ws.Range("B11").Value = "ABS (unpublished) Causes of Death, Australia."

$ws.Range("A12").Value = "References"
$ws.Range("B12").Value = "Council of Australian Governments (COAG), 2011, National Healthcare Agreement."

# B12 picks up its own (slightly larger) font size, which is what drives the
# new style id seen in the saved workbook.
$ws.Range("B12").Font.Size = 12

# Row heights for the touched/new rows, matching the refreshed layout.
$ws.Rows.Item(11).RowHeight = 13.8
$ws.Rows.Item(12).RowHeight = 15

# Leave the selection on the newly added cell, as in the authored edit.
[void]$ws.Range("B12").Select()
